$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.529.95"
$ws.Range("E2").Value = "  -1.52%  "
$ws.Range("D3").Value = "2.901.19"
$ws.Range("E3").Value = "  -2.30%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "526.10"
$ws.Range("E5").Value = "  -2.34%  "
$ws.Range("D6").Value = "142.32"
$ws.Range("E6").Value = "  -5.13%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -3.06%  "
$ws.Range("D9").Value = "2.907.49"
$ws.Range("E9").Value = "  -2.30%  "
$ws.Range("E10").Value = "  -5.23%  "
$ws.Range("D11").Value = "5.96"
$ws.Range("E11").Value = "  -2.66%  "
$ws.Range("D12").Value = "0.359"
$ws.Range("E12").Value = "  -2.26%  "
$ws.Range("D13").Value = "3.410.43"
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("E14").Value = "  +2.26%  "
$ws.Range("D15").Value = "60.532.70"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("D16").Value = "22.64"
$ws.Range("E16").Value = "  -4.02%  "
$ws.Range("D17").Value = "2.904.50"
$ws.Range("E17").Value = "  -2.18%  "
$ws.Range("E18").Value = "  -3.78%  "
$ws.Range("D19").Value = "4.97"
$ws.Range("E19").Value = "  -3.51%  "
$ws.Range("D20").Value = "11.62"
$ws.Range("E20").Value = "  -3.12%  "
$ws.Range("D21").Value = "352.22"
$ws.Range("E21").Value = "  -7.33%  "
$ws.Range("D22").Value = "6.57"
$ws.Range("E22").Value = "  -1.31%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "5.70"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("D25").Value = "64.60"
$ws.Range("E25").Value = "  -1.34%  "
$ws.Range("E26").Value = "  -3.48%  "
$ws.Range("D27").Value = "0.178"
$ws.Range("E27").Value = "  -5.10%  "
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").Value = "7.83"
$ws.Range("E29").Value = "  -4.43%  "
$ws.Range("E30").Value = "  -10.05%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "1.67"
$ws.Range("E32").Value = "  -2.20%  "
$ws.Range("D33").Value = "19.58"
$ws.Range("E33").Value = "  -3.97%  "
$ws.Range("D34").Value = "150.32"
$ws.Range("E34").Value = "  -6.50%  "
$ws.Range("D35").Value = "4.32"
$ws.Range("E35").Value = "  -7.01%  "
$ws.Range("D36").Value = "5.58"
$ws.Range("E36").Value = "  -5.54%  "
$ws.Range("E37").Value = "  -6.46%  "
$ws.Range("E38").Value = "  -5.11%  "
$ws.Range("D39").Value = "37.70"
$ws.Range("E39").Value = "  +0.48%  "
$ws.Range("D40").Value = "1.47"
$ws.Range("E40").Value = "  -5.01%  "
$ws.Range("D41").Value = "3.71"
$ws.Range("E41").Value = "  -4.90%  "
$ws.Range("D42").Value = "2.287.85"
$ws.Range("E42").Value = "  -5.02%  "
$ws.Range("E43").Value = "  -2.92%  "
$ws.Range("D44").Value = "0.0582"
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("D45").Value = "20.49"
$ws.Range("E45").Value = "  -7.40%  "
$ws.Range("D46").Value = "0.997"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "4.97"
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("D48").Value = "0.0237"
$ws.Range("E48").Value = "  -3.27%  "
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("D50").Value = "0.0921"
$ws.Range("E50").Value = "  -3.03%  "
$ws.Range("D51").Value = "248.54"
$ws.Range("E51").Value = "  -6.76%  "
